$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Andrew Case's timecard was updated: 17h 30m -> 18h 30m
$ws.Range("B4").Value = "18h 30m"

# Update the selection to B4, matching the last-saved cursor position
$ws.Range("B4").Select()
